$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'75.827.18"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.99%  '

$ws.Range("D3").Value = "'2.897.00"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.59%  '

$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").Value = "'198.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.27%  '

$ws.Range("D6").Value = "'596.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = "'0.547"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.65%  '

$ws.Range("D9").Value = "'0.198"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.71%  '

$ws.Range("D10").Value = "'2.896.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.55%  '

$ws.Range("D11").Value = "'0.425"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +14.80%  '

$ws.Range("E12").Value = '  -0.09%  '

$ws.Range("D13").Value = "'4.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.37%  '

$ws.Range("D14").Value = "'3.430.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.56%  '

$ws.Range("D15").Value = "'75.788.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.07%  '

$ws.Range("D16").Value = "'0.0000191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.21%  '

$ws.Range("D17").Value = "'27.39"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.11%  '

$ws.Range("D18").Value = "'2.896.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.68%  '

$ws.Range("D19").Value = "'12.70"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.08%  '

$ws.Range("D20").Value = "'8.69"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.40%  '

$ws.Range("D21").Value = "'371.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.23%  '

$ws.Range("D22").Value = "'2.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.93%  '

$ws.Range("D23").Value = "'4.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.46%  '

$ws.Range("B24").Value = 'Dai'
$ws.Range("C24").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.07%  '

$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = "'70.96"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.24%  '

$ws.Range("D26").Value = "'3.061.72"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.03%  '

$ws.Range("D27").Value = "'4.17"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.38%  '

$ws.Range("D28").Value = "'9.62"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.28%  '

$ws.Range("D29").Value = "'0.0000107"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.59%  '

$ws.Range("D30").Value = "'0.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.29%  '

$ws.Range("D31").Value = "'1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.59%  '

$ws.Range("D32").Value = "'498.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -3.14%  '

$ws.Range("D33").Value = "'7.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("D34").Value = "'1.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.96%  '

$ws.Range("D36").Value = "'165.51"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.44%  '

$ws.Range("D37").Value = "'20.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.31%  '

$ws.Range("D38").Value = "'19.64"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.59%  '

$ws.Range("D39").Value = "'0.112"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.90%  '

$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.06%  '

$ws.Range("B41").Value = 'Cronos'
$ws.Range("C41").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D41").Value = "'0.102"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +17.91%  '

$ws.Range("D42").Value = "'179.76"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.06%  '

$ws.Range("D43").Value = "'0.344"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.97%  '

$ws.Range("D44").Value = "'4.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.71%  '

$ws.Range("D45").Value = "'1.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.55%  '

$ws.Range("D46").Value = "'40.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.44%  '

$ws.Range("D47").Value = "'1.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.72%  '

$ws.Range("D48").Value = "'2.31"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.30%  '

$ws.Range("D49").Value = "'0.569"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '

$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").Value = "'3.70"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.15%  '

$ws.Range("B51").Value = 'Mantle'
$ws.Range("C51").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D51").Value = "'0.653"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.06%  '

